$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 2.7
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 3.4
$ws.Range("U2").Value = 1.91
$ws.Range("V2").Value = 1.91
$ws.Range("X2").Value = 13
$ws.Range("Y2").Value = 11
$ws.Range("AA2").Value = 23
$ws.Range("AG2").Value = 301
$ws.Range("AH2").Value = 8
$ws.Range("AK2").Value = 26
$ws.Range("AW2").Value = 4.5
$ws.Range("G12").Value = 2.8
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 2.25
$ws.Range("L12").Value = 2.88
$ws.Range("N12").Value = 19
$ws.Range("Q12").Value = 1.53
$ws.Range("R12").Value = 2.4
$ws.Range("S12").Value = 1.29
$ws.Range("T12").Value = 3.5
$ws.Range("X12").Value = 17
$ws.Range("AA12").Value = 19
$ws.Range("AC12").Value = 19
$ws.Range("AD12").Value = 7.5
$ws.Range("AF12").Value = 34
$ws.Range("AG12").Value = 101
$ws.Range("AH12").Value = 12
$ws.Range("AK12").Value = 23
$ws.Range("AM12").Value = 21
$ws.Range("AN12").Value = 5
$ws.Range("AT12").Value = 3.5
$ws.Range("AX12").Value = 12
$ws.Range("AZ12").Value = 41
$ws.Range("BA12").Value = 51
$ws.Range("BB12").Value = 101
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 8
